$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.848.09"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "'3.119.66"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'577.67"

$ws.Range("D6").Value = "'171.97"
$ws.Range("E6").Value = "  +2.28%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("D9").Value = "'6.37"
$ws.Range("E9").Value = "  -3.47%  "

$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").Value = "'37.11"
$ws.Range("E13").Value = "  +1.97%  "

$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("D15").Value = "'3.635.67"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("D16").Value = "'66.850.87"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").Value = "'3.117.96"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").Value = "'474.35"
$ws.Range("E20").Value = "  +1.69%  "

$ws.Range("E21").Value = "  -0.62%  "

$ws.Range("E22").Value = "  +4.92%  "

$ws.Range("D23").Value = "'83.70"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").Value = "'13.24"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("E25").Value = "  -3.61%  "

$ws.Range("D26").Value = "'10.26"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = "  -1.46%  "

$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("D32").Value = "'0.116"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("E33").Value = "  -6.61%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").Value = "'0.973"

$ws.Range("D37").Value = "'47.10"
$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").Value = "'50.18"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "'2.05"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("E40").Value = "  -2.15%  "

$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("D42").Value = "'8.60"
$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").Value = "'2.812.39"
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").Value = "'382.23"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D47").Value = "'135.63"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D49").Value = "'24.92"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("E51").Value = "  -0.90%  "
